# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Metadata"): update Version, Date, and insert a Jurisdiction row ---
$ws1 = $wb.Worksheets.Item(1)

# Update Version value (row 3)
$ws1.Range("B3").Value2 = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$ws1.Range("B8").Value2 = "2025-10-29T22:15:57+01:00"

# Insert a new row 11 ("Jurisdiction" / blank), pushing subsequent rows down.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("A11").Value2 = "Jurisdiction"
$ws1.Range("B11").Value2 = ""

# --- Sheet 2 ("Elements"): add the II-1 invariant to Performer2.typeId's Constraint(s) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("AJ5").Value2 = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
